$d = $word.ActiveDocument

# 1. Remove the "/{type}" suffix that followed "upload" in the upload-endpoint URL
#    paragraph, leaving just "...rest/upload".
$d.Content.Find.Execute("upload/{type}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "upload", 2) | Out-Null

# 2. Remove the three paragraphs describing "Type can take values: upload (...)",
#    "update (...)" and "unification (...)" in their entirety.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -match "^Type can take values: upload") {
        $startPara = $d.Paragraphs.Item($i)
    }
    if ($text -match "^\s*unification") {
        $endPara = $d.Paragraphs.Item($i)
        break
    }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}

# 3. Relocate the "_GoBack" bookmark from the end of the final paragraph to its
#    start (right after the paragraph mark properties, before the first run).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmStart = $lastPara.Range.Start
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
